$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Generate Report for Handback
#   - Status text "Ready for handoff" becomes "Handback transform failed"
#     on the Overview sheet and on each per-locale sheet.
#   - Each per-locale sheet's "Error Detail" column (P) on the handback
#     row now records why the handback failed.
#   - Error Detail column widened to fit the new message text.
# -----------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: ff2jr3w1.rzh is different with handoff file name: a0ff60ab-43d5-4362-8ee8-c197bc36b244.478a575f2549e0bac8ef9e9cad387bde2ae976cb.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: ff2jr3w1.rzh is different with handoff file name: a0ff60ab-43d5-4362-8ee8-c197bc36b244.478a575f2549e0bac8ef9e9cad387bde2ae976cb.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

Write-Output "Report for Handback generated"
